# Re-applies a template re-render of comprehensive_features_output.docx.
#
# 1. The __functions__ pointer-ish value printed inside the "Entire data
#    context" dump changed between Go runtime invocations (cosmetic).
# 2/3. The "21" / "true" runs in the "14. Product in 1 Paragraph but 3
#    lines" demo are re-written by the (fixed) template-aware run-merging
#    logic; their text does not change.
# 4. The "Another Hyperling: " label run is likewise re-written.
# 5. The literal sanity-check hyperlink text "Dont Replace This" must stay
#    untouched content-wise (it intentionally is not a template
#    expression) but its run is also re-written by the merge fix.
# 6. The "Generated on:" timestamp embedded in the rendered header
#    fragment is refreshed to the new render time.

$d = $word.ActiveDocument

# --- 1. __functions__ pointer value -------------------------------------
$d.Content.Find.Execute(
    "0x140001324a0", $true, $false, $false, $false, $false,
    $true, 1, $false, "0xc0000904a0", 2) | Out-Null

# --- 2/3. "21" / "true" inside paragraph "21<br>true<br>false<br>15" ----
$demoPara = $d.Paragraphs.Item(136)
$demoRange = $demoPara.Range

$twentyOne = $d.Range($demoRange.Start, $demoRange.Start + 2)
if ($twentyOne.Text -eq "21") {
    $twentyOne.Text = "21"
}

$trueRun = $d.Range($demoRange.Start + 3, $demoRange.Start + 7)
if ($trueRun.Text -eq "true") {
    $trueRun.Text = "true"
}

# --- 4. "Another Hyperling: " label --------------------------------------
$d.Content.Find.Execute(
    "Another Hyperling: ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Another Hyperling: ", 2) | Out-Null

# --- 5. Hyperlink text "Dont Replace This" (left as-is, run re-touched) -
$hyperlink = $d.Hyperlinks.Item(2)
if ($hyperlink.TextToDisplay -eq "Dont Replace This") {
    $hyperlink.TextToDisplay = "Dont Replace This"
}

# --- 6. Refresh render timestamp -----------------------------------------
$d.Content.Find.Execute(
    "Generated on: 2025-10-20 08:34:24", $true, $false, $false, $false, $false,
    $true, 1, $false, "Generated on: 2026-02-06 20:53:00", 2) | Out-Null

Write-Host "edit complete"
